$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18 ---------------------------------------------------------
# Set number formats first (matching existing column formats exactly)
# so the engine reuses the existing style indices instead of minting
# new ones.
$ws.Range("A18").NumberFormat = "MM/DD/YY"
$ws.Range("B18:C18").NumberFormat = "HH:MM:SS\ AM/PM"

$ws.Cells.Item(18, 1).Value2 = 41968
$ws.Cells.Item(18, 2).Value2 = 0.625
$ws.Cells.Item(18, 3).Value2 = 0.708333333333333
$ws.Cells.Item(18, 4).Formula = "=ROUND(ABS(C18-B18) * 24, 1)"

# --- Row 19 ---------------------------------------------------------
$ws.Range("A19").NumberFormat = "MM/DD/YY"
$ws.Range("B19:C19").NumberFormat = "HH:MM:SS\ AM/PM"

$ws.Cells.Item(19, 1).Value2 = 41969
$ws.Cells.Item(19, 2).Value2 = 0.416666666666667
$ws.Cells.Item(19, 3).Value2 = 0.625
$ws.Cells.Item(19, 4).Formula = "=ROUND(ABS(C19-B19) * 24, 1)"

# --- Selection --------------------------------------------------------
$ws.Range("D19").Select()
